# Adds non-road hydrogen vehicles, adjusts biofuels share, H2 shares to match
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PVTStL")

# Row 4 corresponds to "aircraft": passenger (B) and freight (C) shares
# reduced from 1 (100%) to 0.83 (83%) to account for non-road hydrogen vehicles.
$ws.Range("B4").Value = 0.83
$ws.Range("C4").Value = 0.83
